$wb = $excel.ActiveWorkbook

# Template sheet to clone: the current last sheet "Bus_Makhulu_r".
$src = $wb.Worksheets.Item("Bus_Makhulu_r")

# --- Add "Truck_Amandla_A2" as a copy of Bus_Makhulu_r, placed at the end ---
$src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$sheetA2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheetA2.Name = "Truck_Amandla_A2"
$sheetA2.Range("H3").Value = "Gear1DShafts1D_Truck_Amandla_A2"

# --- Add "Truck_Amandla_A3" as a copy of Bus_Makhulu_r, placed at the end ---
$src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$sheetA3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheetA3.Name = "Truck_Amandla_A3"
$sheetA3.Range("H3").Value = "Gear1DShafts1D_Truck_Amandla_A3"

# The newly-added last sheet becomes the active / selected tab (mirrors the
# previous state where "Bus_Makhulu_r" - then the last sheet - was active).
$sheetA3.Activate()
